# Enhance start_row/column behaviour on XLS/XLSX
# Shift the whole table down by 3 rows and right by 2 columns (now starts
# at C4 instead of A1), turn the literal boolean cells into TRUE()/FALSE()
# formulas, fix up the number format used by the former A8 (now C11) cell
# so it shares the "TRUE/FALSE" format already used elsewhere, widen the
# newly-exposed column, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the table and 2 blank columns to its left,
# which shifts every existing cell from A1:H9 down to C4:J11.
$ws.Rows("1:3").Insert()
$ws.Columns("A:B").Insert()

# The two boolean cells (now C5/C6) become live formulas instead of
# literal TRUE/FALSE values.
$ws.Range("C5").Formula = "=TRUE()"
$ws.Range("C6").Formula = "=FALSE()"

# The trailing blank cell in the last row (now C11) reuses the
# "TRUE/FALSE" custom number format instead of its own duplicate format.
$ws.Range("C11").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# Column I (the newly exposed "unicode_column" header) gets a wider column.
$ws.Columns("I").ColumnWidth = 17.96

# Move the active selection onto the new data.
$ws.Range("E7:I9").Select()
